# OpenFOAM benchmarks.xlsx edit script
# - Inserts 7 new rows above the old "Medium" results table (old row 9 -> new row 16),
#   creating room for a new per-process-count raw-timing block (rows 10-15) and a new
#   "Parallel Scaling" column (N) for both speed-up figures in the Medium table.
# - Adds the new "Parallel Scaling" header/shared string and its values.
# - Leaves all previously-existing data intact, just shifted down by 7 rows (which
#   Rows.Insert() does natively, the same way the original author's row-insert in
#   Excel would have).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room: insert 7 rows starting at row 9 (old row 9 becomes row 16, old row 24
#    becomes row 31, old row 38 becomes row 45 - matches the target layout exactly).
$ws.Rows("9:15").Insert()

# 2. New column header for the "Medium" table.
$ws.Range("N7").Value = "Parallel Scaling"

# 3. New raw per-process-count timing rows (no #nodes/mesh-gen time recorded for these,
#    only MPI-process count in C and computation time in J), plus their speed-up values.
$ws.Range("C10").Value = 1
$ws.Range("J10").Value = 76737

$ws.Range("C11").Value = 2
$ws.Range("J11").Value = 45419
$ws.Range("N11").Value = 1.69

$ws.Range("C12").Value = 4
$ws.Range("J12").Value = 27010
$ws.Range("N12").Value = 2.84

$ws.Range("C13").Value = 8
$ws.Range("J13").Value = 13804
$ws.Range("N13").Value = 5.56

$ws.Range("C14").Value = 16
$ws.Range("J14").Value = 7520
$ws.Range("N14").Value = 10.2

$ws.Range("C15").Value = 32
$ws.Range("J15").Value = 4706
$ws.Range("N15").Value = 16.31

# 4. Speed-up values for the rest of the (now-shifted) Medium table rows.
$ws.Range("N16").Value = 18.95
$ws.Range("N17").Value = 33.67
$ws.Range("N18").Value = 45.84
$ws.Range("N19").Value = 57.31
$ws.Range("N20").Value = 64.11
$ws.Range("N21").Value = 70.86
$ws.Range("N22").Value = 74.21
$ws.Range("N23").Value = 76.97
$ws.Range("N25").Value = 72.46
$ws.Range("N26").Value = 53

# 5. Restore view state (active cell / selection) to match the saved workbook.
$ws.Range("N27").Select()
